$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the formatting of the other
# header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill the new column with 0 for every data row (rows 2-10)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
